$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-39 down to 28-40
$ws.Rows("27:27").Insert()

# Populate the new row 27 with the new data record
$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(27, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(27, 3).Value = "Metropolitana"
$ws.Cells.Item(27, 4).Value = 44762
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(27, 5).Value = 13
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100102
$ws.Cells.Item(27, 8).Value = "Cítricos"
$ws.Cells.Item(27, 9).Value = 100102006
$ws.Cells.Item(27, 10).Value = "Pomelo"
$ws.Cells.Item(27, 11).Value = "Start Ruby"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 200
$ws.Cells.Item(27, 14).Value = 8000
$ws.Cells.Item(27, 15).Value = 8000
$ws.Cells.Item(27, 16).Value = 8000
$ws.Cells.Item(27, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(27, 18).Value = "Hijuelas"
$ws.Cells.Item(27, 19).Value = 571
$ws.Cells.Item(27, 20).Value = 14
